$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 51.05762969290213

# MACRO_SCORE column (N) was recalculated for rows 2-5; update the cached value.
for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 14).Value = $newValue
}
